$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 117, shifting existing rows 117-194 down to 118-195
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new record's data
$ws.Cells.Item(117, 1).Value = 5
$ws.Cells.Item(117, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(117, 3).Value = 'Maule'
$ws.Cells.Item(117, 4).Value = 44981
$ws.Cells.Item(117, 5).Value = 7
$ws.Cells.Item(117, 6).Value = 100112030
$ws.Cells.Item(117, 7).Value = 'Poroto granado'
$ws.Cells.Item(117, 8).Value = 'Sin especificar'
$ws.Cells.Item(117, 9).Value = 'Primera'
$ws.Cells.Item(117, 10).Value = 400
$ws.Cells.Item(117, 11).Value = 23000
$ws.Cells.Item(117, 12).Value = 24000
$ws.Cells.Item(117, 13).Value = 23500
$ws.Cells.Item(117, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(117, 15).Value = 'Región del Maule'
$ws.Cells.Item(117, 16).Value = 940
$ws.Cells.Item(117, 17).Value = 25
$ws.Cells.Item(117, 18).Value = 'Hortaliza'

# Ensure the date column keeps the same date/time number format used by the
# rest of column D.
$ws.Cells.Item(117, 4).NumberFormat = $ws.Cells.Item(118, 4).NumberFormat
